$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("C4").Value = "Circuitos Elétricos 2 - MCT-2A"
$ws.Range("D4").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "Circuitos Elétricos 2 - ELT-2A"
